$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift existing row1 data from C,E,F to C,D,E by removing the empty column D
$ws.Columns.Item(4).Delete()

# 2. Update the timestamp in E1 to the refreshed value
$ws.Cells.Item(1,5).Value = 45434.701642164349

# 3. Append a new row with the same product info + a newer price check timestamp
$ws.Cells.Item(2,1).Value = $ws.Cells.Item(1,1).Value2
$ws.Cells.Item(2,2).Value = $ws.Cells.Item(1,2).Value2
$ws.Cells.Item(2,3).Value = 3599.97
$ws.Cells.Item(2,4).Value = 3599.97
$ws.Cells.Item(1,5).Copy()
$ws.Cells.Item(2,5).PasteSpecial(-4122)
$ws.Cells.Item(2,5).Value = 45434.702188650903

# 4. Size the new numeric columns
$ws.Columns.Item(3).ColumnWidth = 7.1
$ws.Columns.Item(4).ColumnWidth = 7.1

# 5. Leave the selection on D2, matching the last cell touched interactively
$ws.Range("D2").Select()

$wb.Save()
Write-Host "done"
